# Apply updated crypto price/volume data to match target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.047.78"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "'2.152.64"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'253.22"
$ws.Range("E5").Value = "  +6.06%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'73.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "'39.55"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'2.480.21"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "'2.182.95"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "'41.922.96"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'9.55"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'225.50"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +5.78%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "'10.42"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'36.81"
$ws.Range("E30").Value = "  +10.45%  "
$ws.Range("D31").Value = "'168.03"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'19.86"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'0.0796"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("D38").Value = "'0.0328"
$ws.Range("E38").Value = "  +6.28%  "
$ws.Range("D39").Value = "'12.01"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "'2.04"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("E42").Value = "  -4.68%  "
$ws.Range("D43").Value = "'58.38"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "'99.62"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.23"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D46").Value = "'0.461"
$ws.Range("E46").Value = "  +13.85%  "
$ws.Range("D47").Value = "'0.0960"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'2.37"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +0.79%  "
